# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.497
$ws.Range("C7").Value = -12.666
$ws.Range("A9").Value = -21.79
$ws.Range("C12").Value = -11.555
$ws.Range("A13").Value = -22.025
$ws.Range("C14").Value = -12.369
$ws.Range("E15").Value = 15.782
$ws.Range("A16").Value = -22.016
$ws.Range("A18").Value = -21.956
$ws.Range("C19").Value = -12.059
$ws.Range("A20").Value = -20.242
$ws.Range("A26").Value = -21.61
$ws.Range("C26").Value = -13.169
$ws.Range("A27").Value = -21.898
$ws.Range("C27").Value = -13.25
$ws.Range("E28").Value = 16.84
$ws.Range("A29").Value = -21.315
$ws.Range("C29").Value = -12.007
$ws.Range("E33").Value = 17.303
$ws.Range("A35").Value = -19.982
$ws.Range("E35").Value = 16.604
$ws.Range("A36").Value = -20.391
$ws.Range("C37").Value = -12.854
$ws.Range("C38").Value = -12.652
$ws.Range("E38").Value = 16.667
$ws.Range("E43").Value = 17.023
$ws.Range("E44").Value = 16.67
$ws.Range("A45").Value = -21.617
$ws.Range("E45").Value = 16.695
$ws.Range("C47").Value = -12.331
$ws.Range("E47").Value = 16.184
$ws.Range("C51").Value = -11.432
$ws.Range("E51").Value = 16.812
$ws.Range("C52").Value = -11.644
$ws.Range("E54").Value = 16.536
$ws.Range("A55").Value = -22.06
$ws.Range("C55").Value = -12.929
$ws.Range("A57").Value = -22.22
$ws.Range("E57").Value = 16.227
$ws.Range("E62").Value = 16.376
$ws.Range("E63").Value = 17.585
$ws.Range("E67").Value = 17.341
$ws.Range("A69").Value = -21.699
$ws.Range("C69").Value = -11.312
$ws.Range("C70").Value = -12.065
$ws.Range("E70").Value = 17.354
$ws.Range("A76").Value = -19.999
$ws.Range("C76").Value = -13.076
$ws.Range("A78").Value = -19.854
$ws.Range("C81").Value = -13.605
$ws.Range("E81").Value = 16.838
$ws.Range("A82").Value = -21.976
$ws.Range("A83").Value = -21.768
$ws.Range("C83").Value = -13.075
$ws.Range("E88").Value = 16.218
$ws.Range("A93").Value = -21.435
$ws.Range("C94").Value = -10.913
$ws.Range("E96").Value = 16.196
$ws.Range("A97").Value = -22.132
$ws.Range("E99").Value = 16.738
$ws.Range("C100").Value = -12.869
$ws.Range("C102").Value = -12.922
